$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(80, 8).Value = 1587.28
$ws.Cells.Item(80, 9).Value = 1596.45
$ws.Cells.Item(80, 10).Value = 1550.6
$ws.Cells.Item(80, 11).Value = 4789.35
$ws.Cells.Item(80, 12).Value = 4651.799999999999
$ws.Cells.Item(80, 13).Value = -3791.35
$ws.Cells.Item(80, 14).Value = -6647.799999999999
$ws.Cells.Item(83, 8).Value = 1587.28
$ws.Cells.Item(83, 9).Value = 1596.45
$ws.Cells.Item(83, 10).Value = 1550.6
$ws.Cells.Item(83, 11).Value = 14368.05
$ws.Cells.Item(83, 12).Value = 13955.4
$ws.Cells.Item(83, 13).Value = -9376.050000000001
$ws.Cells.Item(83, 14).Value = -23939.4
$ws.Cells.Item(98, 8).Value = 2349.5454
$ws.Cells.Item(98, 9).Value = 2373.724
$ws.Cells.Item(98, 10).Value = 2174.25
$ws.Cells.Item(98, 11).Value = 2373.724
$ws.Cells.Item(98, 12).Value = 2174.25
$ws.Cells.Item(98, 13).Value = -875.7240000000002
$ws.Cells.Item(98, 14).Value = -5170.25
$ws.Cells.Item(112, 8).Value = 1374.8148
$ws.Cells.Item(112, 10).Value = 1385.3846
$ws.Cells.Item(112, 12).Value = 4156.1538
$ws.Cells.Item(112, 14).Value = -6372.1538
$ws.Cells.Item(122, 8).Value = 2349.5454
$ws.Cells.Item(122, 9).Value = 2373.724
$ws.Cells.Item(122, 10).Value = 2174.25
$ws.Cells.Item(122, 11).Value = 7121.172
$ws.Cells.Item(122, 12).Value = 6522.75
$ws.Cells.Item(122, 13).Value = -4671.172
$ws.Cells.Item(122, 14).Value = -11422.75
$ws.Cells.Item(127, 8).Value = 3354.2856
$ws.Cells.Item(127, 9).Value = 2540
$ws.Cells.Item(127, 11).Value = 7620
$ws.Cells.Item(127, 13).Value = -2660
$ws.Cells.Item(137, 8).Value = 1075.0652
$ws.Cells.Item(137, 9).Value = 898.5714
$ws.Cells.Item(137, 10).Value = 1349.6111
$ws.Cells.Item(137, 11).Value = 2695.7142
$ws.Cells.Item(137, 12).Value = 4048.8333
$ws.Cells.Item(137, 13).Value = -145.7142000000003
$ws.Cells.Item(137, 14).Value = -9148.8333
$ws.Cells.Item(138, 8).Value = 1852.2344
$ws.Cells.Item(138, 9).Value = 1590.48
$ws.Cells.Item(138, 10).Value = 2787.0715
$ws.Cells.Item(138, 11).Value = 4771.440000000001
$ws.Cells.Item(138, 12).Value = 8361.2145
$ws.Cells.Item(138, 13).Value = 368.5599999999995
$ws.Cells.Item(138, 14).Value = -18641.2145
$ws.Cells.Item(140, 8).Value = 62993.316
$ws.Cells.Item(140, 10).Value = 62993.316
$ws.Cells.Item(140, 12).Value = 62993.316
$ws.Cells.Item(140, 14).Value = -73353.31599999999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2900.3933
$ws.Cells.Item(32, 9).Value = 2437
$ws.Cells.Item(32, 11).Value = 2437
$ws.Cells.Item(32, 13).Value = -2150
$ws.Cells.Item(44, 8).Value = 29983.334
$ws.Cells.Item(44, 10).Value = 29983.334
$ws.Cells.Item(44, 12).Value = 29983.334
$ws.Cells.Item(44, 14).Value = -30959.334
$ws.Cells.Item(55, 8).Value = 17140
$ws.Cells.Item(55, 10).Value = 17140
$ws.Cells.Item(55, 12).Value = 17140
$ws.Cells.Item(55, 14).Value = -17770
$ws.Cells.Item(61, 8).Value = 2147.9697
$ws.Cells.Item(61, 9).Value = 1576.871
$ws.Cells.Item(61, 11).Value = 1576.871
$ws.Cells.Item(61, 13).Value = -1364.871
$ws.Cells.Item(80, 8).Value = 24550
$ws.Cells.Item(80, 9).Value = 100
$ws.Cells.Item(80, 10).Value = 49000
$ws.Cells.Item(80, 11).Value = 100
$ws.Cells.Item(80, 12).Value = 49000
$ws.Cells.Item(80, 13).Value = 898
$ws.Cells.Item(80, 14).Value = -50996
$ws.Cells.Item(83, 8).Value = 24550
$ws.Cells.Item(83, 9).Value = 100
$ws.Cells.Item(83, 10).Value = 49000
$ws.Cells.Item(83, 11).Value = 300
$ws.Cells.Item(83, 12).Value = 147000
$ws.Cells.Item(83, 13).Value = 4692
$ws.Cells.Item(83, 14).Value = -156984
$ws.Cells.Item(102, 8).Value = 1399.8572
$ws.Cells.Item(102, 9).Value = 1399.8572
$ws.Cells.Item(102, 11).Value = 1399.8572
$ws.Cells.Item(102, 13).Value = 222.1428000000001
$ws.Cells.Item(122, 8).Value = 1527
$ws.Cells.Item(122, 9).Value = 1399.5652
$ws.Cells.Item(122, 11).Value = 4198.6956
$ws.Cells.Item(122, 13).Value = -1748.6956
$ws.Cells.Item(132, 8).Value = 1145.3412
$ws.Cells.Item(132, 9).Value = 852.3280999999999
$ws.Cells.Item(132, 10).Value = 2038.3334
$ws.Cells.Item(132, 11).Value = 2556.9843
$ws.Cells.Item(132, 12).Value = 6115.0002
$ws.Cells.Item(132, 13).Value = -26.98430000000008
$ws.Cells.Item(132, 14).Value = -11175.0002
$ws.Cells.Item(136, 8).Value = 2147.9697
$ws.Cells.Item(136, 9).Value = 1576.871
$ws.Cells.Item(136, 11).Value = 4730.613
$ws.Cells.Item(136, 13).Value = -2180.613

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(82, 8).Value = 36500
$ws.Cells.Item(82, 9).Value = 25000
$ws.Cells.Item(82, 11).Value = 25000
$ws.Cells.Item(82, 13).Value = -24617
$ws.Cells.Item(85, 8).Value = 36500
$ws.Cells.Item(85, 9).Value = 25000
$ws.Cells.Item(85, 11).Value = 25000
$ws.Cells.Item(85, 13).Value = -23674
$ws.Cells.Item(134, 8).Value = 3321.739
$ws.Cells.Item(134, 9).Value = 3283.7727
$ws.Cells.Item(134, 10).Value = 4157
$ws.Cells.Item(134, 11).Value = 9851.3181
$ws.Cells.Item(134, 12).Value = 12471
$ws.Cells.Item(134, 13).Value = -7316.3181
$ws.Cells.Item(134, 14).Value = -17541

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2443.6316
$ws.Cells.Item(31, 9).Value = 2237.2222
$ws.Cells.Item(31, 10).Value = 2629.4
$ws.Cells.Item(31, 11).Value = 2237.2222
$ws.Cells.Item(31, 12).Value = 2629.4
$ws.Cells.Item(31, 13).Value = -1942.2222
$ws.Cells.Item(31, 14).Value = -3219.4
$ws.Cells.Item(34, 8).Value = 2443.6316
$ws.Cells.Item(34, 9).Value = 2237.2222
$ws.Cells.Item(34, 10).Value = 2629.4
$ws.Cells.Item(34, 11).Value = 2237.2222
$ws.Cells.Item(34, 12).Value = 2629.4
$ws.Cells.Item(34, 13).Value = -2035.2222
$ws.Cells.Item(34, 14).Value = -3033.4
$ws.Cells.Item(41, 8).Value = 29000
$ws.Cells.Item(41, 10).Value = 29000
$ws.Cells.Item(41, 12).Value = 29000
$ws.Cells.Item(41, 14).Value = -29856
$ws.Cells.Item(50, 8).Value = 14426.667
$ws.Cells.Item(50, 10).Value = 14426.667
$ws.Cells.Item(50, 12).Value = 14426.667
$ws.Cells.Item(50, 14).Value = -15676.667
$ws.Cells.Item(51, 8).Value = 35000
$ws.Cells.Item(51, 10).Value = 35000
$ws.Cells.Item(51, 12).Value = 35000
$ws.Cells.Item(51, 14).Value = -36472
$ws.Cells.Item(58, 8).Value = 1061612.1
$ws.Cells.Item(58, 9).Value = 1499992.8
$ws.Cells.Item(58, 10).Value = 2192.25
$ws.Cells.Item(58, 11).Value = 1499992.8
$ws.Cells.Item(58, 12).Value = 2192.25
$ws.Cells.Item(58, 13).Value = -1499789.8
$ws.Cells.Item(58, 14).Value = -2598.25
$ws.Cells.Item(59, 8).Value = 17800
$ws.Cells.Item(59, 10).Value = 17800
$ws.Cells.Item(59, 12).Value = 17800
$ws.Cells.Item(59, 14).Value = -20090
$ws.Cells.Item(60, 8).Value = 20472.867
$ws.Cells.Item(60, 10).Value = 20472.867
$ws.Cells.Item(60, 12).Value = 20472.867
$ws.Cells.Item(60, 14).Value = -21494.867
$ws.Cells.Item(61, 8).Value = 35000
$ws.Cells.Item(61, 10).Value = 35000
$ws.Cells.Item(61, 12).Value = 35000
$ws.Cells.Item(61, 14).Value = -35696
$ws.Cells.Item(122, 8).Value = 4165
$ws.Cells.Item(122, 9).Value = 3639.1428
$ws.Cells.Item(122, 10).Value = 4533.1
$ws.Cells.Item(122, 11).Value = 10917.4284
$ws.Cells.Item(122, 12).Value = 13599.3
$ws.Cells.Item(122, 13).Value = -8467.428400000001
$ws.Cells.Item(122, 14).Value = -18499.3
$ws.Cells.Item(132, 8).Value = 1313
$ws.Cells.Item(132, 9).Value = 1070.3914
$ws.Cells.Item(132, 11).Value = 3211.1742
$ws.Cells.Item(132, 13).Value = -681.1741999999999
$ws.Cells.Item(134, 8).Value = 1463.4237
$ws.Cells.Item(134, 9).Value = 1350.9783
$ws.Cells.Item(134, 11).Value = 4052.9349
$ws.Cells.Item(134, 13).Value = -1517.9349
$ws.Cells.Item(136, 8).Value = 1061612.1
$ws.Cells.Item(136, 9).Value = 1499992.8
$ws.Cells.Item(136, 10).Value = 2192.25
$ws.Cells.Item(136, 11).Value = 4499978.4
$ws.Cells.Item(136, 12).Value = 6576.75
$ws.Cells.Item(136, 13).Value = -4497428.4
$ws.Cells.Item(136, 14).Value = -11676.75

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(107, 8).Value = 466.3158
$ws.Cells.Item(107, 10).Value = 462.9375
$ws.Cells.Item(107, 12).Value = 1388.8125
$ws.Cells.Item(107, 14).Value = -5228.8125
$ws.Cells.Item(121, 8).Value = 655.4
$ws.Cells.Item(121, 10).Value = 713.4
$ws.Cells.Item(121, 12).Value = 2140.2
$ws.Cells.Item(121, 14).Value = -4760.2
$ws.Cells.Item(132, 8).Value = 943.5
$ws.Cells.Item(132, 9).Value = 849.6667
$ws.Cells.Item(132, 11).Value = 7647.0003
$ws.Cells.Item(132, 13).Value = -5117.0003
$ws.Cells.Item(138, 8).Value = 2339
$ws.Cells.Item(138, 9).Value = 2082
$ws.Cells.Item(138, 11).Value = 6246
$ws.Cells.Item(138, 13).Value = -1106

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 5000
$ws.Cells.Item(80, 10).Value = 5000
$ws.Cells.Item(80, 12).Value = 5000
$ws.Cells.Item(80, 14).Value = -6996
$ws.Cells.Item(83, 8).Value = 5000
$ws.Cells.Item(83, 10).Value = 5000
$ws.Cells.Item(83, 12).Value = 25000
$ws.Cells.Item(83, 14).Value = -34984
$ws.Cells.Item(102, 8).Value = 2316.2964
$ws.Cells.Item(102, 10).Value = 2411.7144
$ws.Cells.Item(102, 12).Value = 2411.7144
$ws.Cells.Item(102, 14).Value = -5655.7144
$ws.Cells.Item(122, 8).Value = 1563.25
$ws.Cells.Item(122, 9).Value = 1516.2142
$ws.Cells.Item(122, 10).Value = 1673
$ws.Cells.Item(122, 11).Value = 4548.642599999999
$ws.Cells.Item(122, 12).Value = 5019
$ws.Cells.Item(122, 13).Value = -2098.642599999999
$ws.Cells.Item(122, 14).Value = -9919
$ws.Cells.Item(132, 8).Value = 621832.7
$ws.Cells.Item(132, 9).Value = 786078
$ws.Cells.Item(132, 10).Value = 2754.2307
$ws.Cells.Item(132, 11).Value = 2358234
$ws.Cells.Item(132, 12).Value = 8262.6921
$ws.Cells.Item(132, 13).Value = -2355704
$ws.Cells.Item(132, 14).Value = -13322.6921
$ws.Cells.Item(140, 8).Value = 44127.3
$ws.Cells.Item(140, 10).Value = 44127.3
$ws.Cells.Item(140, 12).Value = 44127.3
$ws.Cells.Item(140, 14).Value = -54487.3

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 2566.842
$ws.Cells.Item(61, 9).Value = 2439.5
$ws.Cells.Item(61, 10).Value = 2785.1428
$ws.Cells.Item(61, 11).Value = 2439.5
$ws.Cells.Item(61, 12).Value = 2785.1428
$ws.Cells.Item(61, 13).Value = -2237.5
$ws.Cells.Item(61, 14).Value = -3189.1428
$ws.Cells.Item(113, 8).Value = 2566.842
$ws.Cells.Item(113, 9).Value = 2439.5
$ws.Cells.Item(113, 10).Value = 2785.1428
$ws.Cells.Item(113, 11).Value = 2439.5
$ws.Cells.Item(113, 12).Value = 2785.1428
$ws.Cells.Item(113, 13).Value = -269.5
$ws.Cells.Item(113, 14).Value = -7125.1428

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(38, 8).Value = 32352.334
$ws.Cells.Item(38, 10).Value = 32352.334
$ws.Cells.Item(38, 12).Value = 32352.334
$ws.Cells.Item(38, 14).Value = -33298.334
$ws.Cells.Item(107, 8).Value = 614.8
$ws.Cells.Item(107, 9).Value = 317.81818
$ws.Cells.Item(107, 10).Value = 977.7778
$ws.Cells.Item(107, 11).Value = 953.45454
$ws.Cells.Item(107, 12).Value = 2933.3334
$ws.Cells.Item(107, 13).Value = 966.54546
$ws.Cells.Item(107, 14).Value = -6773.3334
$ws.Cells.Item(136, 8).Value = 12079509
$ws.Cells.Item(136, 9).Value = 18521170
$ws.Cells.Item(136, 11).Value = 55563510
$ws.Cells.Item(136, 13).Value = -55560960
